$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shift-assignment values for rows 3..19 (columns D, E, F)
$data = @{
    3  = @("Kapil",   "Divik",   "Sushvin")
    4  = @("Sushvin",  "Kapil",   "Naveen")
    5  = @("Divik",    "Naveen",  "Kapil")
    6  = @("Kapil",    "Sushvin", "Divik")
    7  = @("Sushvin",  "Divik",   "Naveen")
    8  = @("Naveen",   "Kapil",   "Sushvin")
    9  = @("Sushvin",  "Divik",   "Naveen")
    10 = @("Divik",    "Naveen",  "Kapil")
    11 = @("Naveen",   "Kapil",   "Sushvin")
    12 = @("Sushvin",  "Divik",   "Naveen")
    13 = @("Kapil",    "Naveen",  "Divik")
    14 = @("Naveen",   "Sushvin", "Kapil")
    15 = @("Kapil",    "Divik",   "Sushvin")
    16 = @("Divik",    "Sushvin", "Kapil")
    17 = @("Naveen",   "Kapil",   "Divik")
    18 = @("Sushvin",  "Naveen",  "Kapil")
    19 = @("Divik",    "Kapil",   "Naveen")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
    $ws.Range("F$row").Value = $vals[2]
}

# Remove the now-obsolete trailing rows (20, 21, 22)
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(20).Delete()

# Update the date-time number formats to date-only formats
foreach ($cell in $ws.UsedRange.Cells) {
    if ($cell.NumberFormat -eq "yyyy-mm-dd h:mm:ss") {
        $cell.NumberFormat = "yyyy-mm-dd"
    } elseif ($cell.NumberFormat -eq "YYYY-MM-DD HH:MM:SS") {
        $cell.NumberFormat = "YYYY-MM-DD"
    }
}
